# Slide 24 ("Loop and Exit Statements (continued)"), Content Placeholder 2:
# the first bullet currently reads
#   "CPRL also permits an optional for prefix for a loop."
# and needs to become three runs so that "for" is wrapped in curly quotes
# and rendered in the Consolas (code) font, matching the other code
# snippets used throughout the deck:
#   "CPRL also permits an optional “" + "for" (Consolas) + "” prefix for a loop."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$oldFirstLine = "CPRL also permits an optional for prefix for a loop."
$newFirstLine = "CPRL also permits an optional " + $openQuote + "for" + $closeQuote + " prefix for a loop."

# Replace just the first paragraph's text (it has no trailing paragraph
# mark within this length, so Characters(1, N) stays inside paragraph 1).
$paraRange = $tr.Characters(1, $oldFirstLine.Length)
$paraRange.Text = $newFirstLine

# Re-apply the Consolas font to the "for" that is now wrapped in quotes.
$forStart = $newFirstLine.IndexOf($openQuote) + 2
$forRange = $tr.Characters($forStart, 3)
$forRange.Font.Name = "Consolas"
